# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# Commit: "Updated cryptos list on Thu Jul 18 08:49:29 UTC 2024 with GitHub Actions"
#
# Price (col D) and 1h-volume (col E) text is refreshed for most rows; two
# pairs of rows (23/24 and 42/43) had their coin identity swapped in the feed
# (Polygon<->Dai, OKB<->Filecoin), so Coin/Link/Price/Volume are all rewritten
# for those four rows.
#
# Every value in this sheet is stored as literal text (the source feed already
# renders formatted strings like "64.759.48" or "  -0.81%  "). Some Price cells
# look like plain decimals (e.g. "571.95"), and Excel would normally infer those
# as numbers on a bare .Value assignment, silently reformatting them (trailing
# zeros dropped, binary float noise, etc). Marking the cell as Text (@) before
# the assignment keeps the exact string, then ClearFormats() drops the helper
# number-format again so the cell is left with no explicit style, matching the
# original (unstyled) data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '64.759.48'
$ws.Range('E2').Value = '  -0.81%  '

# Row 3
$ws.Range('D3').Value = '3.441.60'
$ws.Range('E3').Value = '  -1.27%  '

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.95'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.26%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.43'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.28%  '

# Row 7
$ws.Range('E7').Value = '  +0.01%  '

# Row 8
$ws.Range('D8').Value = '3.441.44'
$ws.Range('E8').Value = '  -1.26%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.573'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -6.95%  '

# Row 10
$ws.Range('E10').Value = '  -1.26%  '

# Row 11
$ws.Range('E11').Value = '  -3.92%  '

# Row 12
$ws.Range('E12').Value = '  -1.62%  '

# Row 13
$ws.Range('D13').Value = '4.031.93'
$ws.Range('E13').Value = '  -1.39%  '

# Row 14
$ws.Range('E14').Value = '  -0.55%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.53'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.25%  '

# Row 16
$ws.Range('E16').Value = '  -9.91%  '

# Row 17
$ws.Range('D17').Value = '64.828.63'
$ws.Range('E17').Value = '  -0.65%  '

# Row 18
$ws.Range('D18').Value = '3.412.98'
$ws.Range('E18').Value = '  -2.35%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.21'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.97%  '

# Row 20
$ws.Range('E20').Value = '  -4.45%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '376.57'
$ws.Range('D21').ClearFormats()

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.95'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.43%  '

# Row 23
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.15%  '

# Row 24
$ws.Range('B24').Value = 'Polygon'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.541'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.78%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '71.85'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.27%  '

# Row 26
$ws.Range('E26').Value = '  -0.94%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.81'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.12%  '

# Row 28
$ws.Range('E28').Value = '  -0.71%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.14%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.45'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.86%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.04'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.31%  '

# Row 32
$ws.Range('E32').Value = '  -2.76%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.07'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.40%  '

# Row 34
$ws.Range('E34').Value = '  -4.35%  '

# Row 35
$ws.Range('E35').Value = '  -1.92%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '160.76'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.43%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.85'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.88%  '

# Row 38
$ws.Range('D38').Value = '2.904.35'
$ws.Range('E38').Value = '  -4.50%  '

# Row 39
$ws.Range('E39').Value = '  -4.12%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.64'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.94%  '

# Row 41
$ws.Range('E41').Value = '  -3.79%  '

# Row 42
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.90'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.07%  '

# Row 43
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.51'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.56%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.780'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.16%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.83'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.74%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0310'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.96%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.27'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.88%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '316.82'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.67%  '

# Row 49
$ws.Range('E49').Value = '  -2.94%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.46'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.57%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.843'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.46%  '
